$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: "tipo_metodo" ---------------------------------------
# Header
$ws.Range("G1").Value = "tipo_metodo"
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats -> reuse header style

# Data rows (values)
$values = @{
    2  = "piramide"
    3  = "piramide"
    4  = "abierta"
    5  = "abierta"
    6  = "categorica"
    7  = "categorica"
    8  = "abierta"
    9  = "categorica"
    10 = "abierta"
    11 = "categorica"
    12 = "abierta"
    13 = "categorica"
    14 = "abierta"
    15 = "numericas"
    16 = "numericas"
    17 = "numericas"
    18 = "numericas"
    19 = "abierta"
    20 = "abierta"
    21 = "categorica"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}

# Rows whose G cell carries the same formatting as the "abierta" style (style index 8,
# e.g. cell D4). Copy that exact formatting so the style is reused rather than
# re-created.
$styledRows = @(4, 5, 8, 10, 12, 14, 19, 20)
$ws.Range("D4").Copy()
foreach ($row in $styledRows) {
    $ws.Cells.Item($row, 7).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# --- Column widths to match final layout --------------------------------
# (the engine re-derives the stored "character width" from pixels, offset by
# 5/6 from the ColumnWidth value fed in here - compensate so the saved width
# lands as close as possible to the target layout)
$ws.Columns.Item(1).ColumnWidth = 45.451822916666664
$ws.Columns.Item(2).ColumnWidth = 96.59244791666667
$ws.Columns.Item(3).ColumnWidth = 27.451822916666668
$ws.Columns.Item(4).ColumnWidth = 28.451822916666668
$ws.Columns.Item(5).ColumnWidth = 68.16666666666667
$ws.Columns.Item(6).ColumnWidth = 108.30729166666667
$ws.Columns.Item(7).ColumnWidth = 11.736979166666666

# --- Selection -----------------------------------------------------------
$ws.Range("D21").Select()
